$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New status notes added in column G, next to the two rows that were
# recently updated (Anthony / Tommy).
$ws.Range("G7").Value = "Tom a fini 23 Oct 2014 / Demande lui ses heures"
$ws.Range("G5").Value = "Commencé par Antho 23 Oct 2014 - 21:30"

# Match the formatting already used by the rest of the data rows
# (centered horizontally/vertically) without introducing new style
# records - copy the format from a sibling cell on the same row.
$ws.Range("A5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("G7").PasteSpecial(-4122)

# Size column G to fit the new text.
$ws.Columns("G").ColumnWidth = 43.14

# Extend the "contains Terminé" highlight rule to the two new cells.
$ws.Range("G5").FormatConditions.Add(2, 0, "Terminé")
$ws.Range("G7").FormatConditions.Add(2, 0, "Terminé")

# Move the active selection to reflect where the user ended up.
$ws.Range("G3").Select()
